$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values that are stored as text in the source
# workbook (some are not valid numbers at all, e.g. "91.281.21"; others like
# "1.00" would silently collapse to the number 1 if Excel auto-typed them).
# Force each Price cell we touch to stay text-formatted before writing it.

$ws.Range("D2").NumberFormat = "@"
$ws.Range('D2').Value = '91.281.21'
$ws.Range('E2').Value = '  +4.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range('D3').Value = '3.120.60'
$ws.Range('E3').Value = '  +2.35%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '219.58'
$ws.Range('E5').Value = '  +4.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '624.04'
$ws.Range('E6').Value = '  +1.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '0.389'
$ws.Range('E7').Value = '  +5.86%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '0.955'
$ws.Range('E8').Value = '  +20.28%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '3.116.68'
$ws.Range('E10').Value = '  +2.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.722'
$ws.Range('E11').Value = '  +21.33%  '

$ws.Range('E12').Value = '  +6.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '0.0000257'
$ws.Range('E13').Value = '  +9.05%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '34.46'
$ws.Range('E14').Value = '  +8.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '91.082.77'
$ws.Range('E15').Value = '  +3.96%  '

$ws.Range('E16').Value = '  +2.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '3.695.74'
$ws.Range('E17').Value = '  +2.23%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '3.129.78'
$ws.Range('E18').Value = '  +2.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '3.78'
$ws.Range('E19').Value = '  +17.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '0.0000221'
$ws.Range('E20').Value = '  +11.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '14.07'
$ws.Range('E21').Value = '  +7.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '435.23'
$ws.Range('E22').Value = '  +4.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '8.75'
$ws.Range('E23').Value = '  +8.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '5.17'
$ws.Range('E24').Value = '  +6.61%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '6.16'
$ws.Range('E25').Value = '  +13.29%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '86.18'
$ws.Range('E26').Value = '  +5.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '12.16'
$ws.Range('E27').Value = '  +4.52%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '3.288.51'
$ws.Range('E28').Value = '  +1.87%  '

$ws.Range('E29').Value = '  +0.03%  '

$ws.Range('E30').Value = '  -2.61%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '9.04'
$ws.Range('E31').Value = '  +13.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -7.71%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '529.35'
$ws.Range('E33').Value = '  +4.72%  '

$ws.Range('E34').Value = '  +6.71%  '

$ws.Range('E35').Value = '  +6.08%  '

$ws.Range('E36').Value = '  +9.63%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '23.55'
$ws.Range('E37').Value = '  +6.49%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '1.86'
$ws.Range('E38').Value = '  +3.96%  '

$ws.Range('E39').Value = '  +3.21%  '

$ws.Range('E40').Value = '  +0.41%  '

$ws.Range('E41').Value = '  -0.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '0.0829'
$ws.Range('E42').Value = '  +23.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '0.148'
$ws.Range('E43').Value = '  +12.82%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '0.378'
$ws.Range('E45').Value = '  +5.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '1.91'
$ws.Range('E46').Value = '  +6.70%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '146.93'
$ws.Range('E47').Value = '  -0.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '43.86'
$ws.Range('E48').Value = '  +1.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '166.38'
$ws.Range('E50').Value = '  +7.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').Value = '0.000258'
$ws.Range('E51').Value = '  +22.34%  '
